# Ajout des choses dans les docs
# This script reproduces a re-organisation of the weekly log sheet:
#  - 3 new "week" header rows (merged, yellow, centered) are inserted
#  - a few blank separator rows are inserted
#  - the text of one existing entry is updated
#  - a brand new entry (week 3) is appended
#  - the sheet is set up for printing on a specific paper size

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1. Make room for the new rows by inserting blank rows, working
#    from the bottom of the sheet upward so row numbers stay valid.
# ---------------------------------------------------------------
$ws.Rows.Item(11).Insert()   # -> becomes blank row 12
$ws.Rows.Item(10).Insert()   # -> becomes blank row 14
$ws.Rows.Item(9).Insert()    # -> becomes blank row 12 (pre "2eme semaine" header slot)
$ws.Rows.Item(7).Insert()    # -> becomes "2eme semaine" header slot (row 9)
$ws.Rows.Item(4).Insert()    # -> becomes blank row 5
$ws.Rows.Item(1).Insert()    # -> becomes "1er semaine" header slot (row 1)

# ---------------------------------------------------------------
# 2. Fix up the value of the entry that is now on row 19 (former
#    row 13) and shrink its row height from 60 to 30.
# ---------------------------------------------------------------
$ws.Range("B19").Value = "Modification du MLD, j'ai ajouté une table pour le paiement des articles, j'ai décidé de laissé la table taille et couleur dans mon MLD"
$ws.Rows.Item(19).RowHeight = 30

# ---------------------------------------------------------------
# 3. Add the three week-separator header rows (merged A:C, yellow
#    fill, centered horizontally).
# ---------------------------------------------------------------
function Add-WeekHeader([int]$row, [string]$text) {
    $rng = $ws.Range("A" + $row + ":C" + $row)
    $rng.Interior.Color = 65535
    $rng.HorizontalAlignment = -4108
    $ws.Range("A" + $row).Value = $text
    $rng.Merge()
}

Add-WeekHeader 1 "1er semaine"
Add-WeekHeader 9 "2ème semaine"

# ---------------------------------------------------------------
# 4. Append the new "3eme semaine" header and its single log entry
#    at the bottom of the sheet (rows 20 and 21).
# ---------------------------------------------------------------
Add-WeekHeader 20 "3ème semaine "

$ws.Range("A21").Value = [DateTime]::new(2018, 2, 13)
$ws.Range("A21").NumberFormat = "m/d/yy"

$ws.Range("B21").Value = "J'ai modifié ma base de donnée en ajoutant un nouveau champ dans la table article.`nJ'ai ajouté des choses dans mes use case, j'ai aussi modifié plusieurs pages du site internet pour que les liens entre les pages soient justes. J'ai supprimé des blocs du site que je n'allais pas utiliser. Ajout d'une nouvelle table dans mon MLD"
$ws.Range("B21").WrapText = $true

$ws.Range("C21").Value = "6 périodes"

$ws.Rows.Item(21).RowHeight = 60

# ---------------------------------------------------------------
# 5. Update the view: scroll down a bit and select C22 (the cell
#    just below the new last row), matching where the user's
#    cursor ended up after the edit.
# ---------------------------------------------------------------
$ws.Range("C22").Select()

# ---------------------------------------------------------------
# 6. Configure the page for printing (paper size 152 = Japanese
#    postcard-ish "B5" class size used by the source workbook,
#    portrait orientation).
# ---------------------------------------------------------------
$ws.PageSetup.PaperSize = 152
$ws.PageSetup.Orientation = 1
